$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Alineado a la pagina" -- shrink/align the letterhead logo so it fits
#    the page precisely. The wp:extent was already correct; the picture's
#    own shape size (pic:spPr/a:xfrm/a:ext) still held the old, slightly
#    larger dimensions. Re-set Width/Height (in points) so both the
#    wrapper extent and the shape's internal transform agree.
# ---------------------------------------------------------------------------
$logo = $d.Shapes(1)
$logo.Width  = 2705682 / 12700.0
$logo.Height = 1409065 / 12700.0

# ---------------------------------------------------------------------------
# 2) Split the heading "CENTRAL DE ACEROS M Y M" in two runs, with the
#    _GoBack bookmark (last-edit marker) sitting right at the split point --
#    this is what Word leaves behind when that text was last touched there.
#    Inserting the (still empty) bookmark at the exact split position, while
#    the original text is still intact, makes the engine break the run in
#    place without disturbing the shared run formatting.
# ---------------------------------------------------------------------------
$heading = $d.Content
$found = $heading.Find.Execute("CENTRAL DE ACEROS M Y M")
if (-not $found) {
    throw "Could not locate the heading text to split."
}
$splitPos = $heading.Start + "CENTRAL DE ACEROS".Length
$splitPoint = $d.Range($splitPos, $splitPos)

# Re-homing the bookmark under its existing name moves it here (and removes
# it from its previous location further down the document) instead of
# creating a duplicate.
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
